# Apply the "changing the sheet name" edit:
#  - Duplicate the current "test" sheet (template) and place the copy at the end,
#    keeping its original content (A1 timestamp stays 2024/08/06 19:46:23).
#  - Rename the old "test1" sheet to "test2" and stamp its A1 with the new time.
#  - Rename the old "test" sheet to "test1" and stamp its A1 with the previous
#    "test1" timestamp.
#  - Rename the freshly duplicated sheet to "test" (it keeps the original values).

$wb = $excel.ActiveWorkbook

$sheetTest1 = $wb.Worksheets.Item("test1")
$sheetTest  = $wb.Worksheets.Item("test")

# 1) Copy the "test" sheet to the end of the workbook before any renames/edits,
#    so the copy retains the original "test" content (A1 = 19:46:23).
$sheetTest.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# 2) Rename old "test1" -> "test2" and update its timestamp.
$sheetTest1.Name = "test2"
$sheetTest1.Range("A1").Value = "2024/08/06 19:58:11"

# 3) Rename old "test" -> "test1" and update its timestamp.
$sheetTest.Name = "test1"
$sheetTest.Range("A1").Value = "2024/08/06 19:46:41"

# 4) Rename the duplicated sheet to "test" (content/timestamp unchanged).
$newSheet.Name = "test"

# Restore the originally active tab (the workbook view itself is untouched by
# this edit, so keep the first sheet selected/active like before).
$sheetTest1.Activate()
